$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292:356 down to 293:357.
# This mirrors Excel's "Insert" on a whole row, which also shifts formatting down
# and moves the sheet's used-range dimension from A1:R356 to A1:R357.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new data record.
$ws.Range("A292").Value = 7
$ws.Range("B292").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C292").Value = "Ñuble"
$ws.Range("D292").Value = 45204
$ws.Range("E292").Value = 16
$ws.Range("F292").Value = 100112032
$ws.Range("G292").Value = "Zapallo italiano"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Primera"
$ws.Range("J292").Value = 100
$ws.Range("K292").Value = 14000
$ws.Range("L292").Value = 14000
$ws.Range("M292").Value = 14000
$ws.Range("N292").Value = "$/caja 50 unidades"
$ws.Range("O292").Value = "Región de Arica y Parinacota"
$ws.Range("P292").Value = 280
$ws.Range("Q292").Value = 50
$ws.Range("R292").Value = "Hortaliza"
